$d = $word.ActiveDocument

# 1. Fix the merge-field date placeholder: add an explicit date-format spec
#    "[res_letterbox.doc_date]" -> "[res_letterbox.doc_date;frm=dd/mm/yyyy]"
# The "_GoBack" bookmark sits right after "doc_date"; temporarily drop it so
# the new text lands on the correct side of the boundary, then restore it
# at the same spot (mirrors Word's own "last edit" bookmark behaviour).
$hadGoBack = $d.Bookmarks.Exists("_GoBack")
if ($hadGoBack) {
    $goBackStart = $d.Bookmarks.Item("_GoBack").Start
    $d.Bookmarks.Item("_GoBack").Delete()
} else {
    $goBackStart = $null
}

$d.Content.Find.Execute("doc_date", $true, $false, $false, $false, $false,
                         $true, 1, $false, "doc_date;frm=dd/mm/yyyy", 2)

if ($hadGoBack) {
    $d.Bookmarks.Add("_GoBack", $d.Range($goBackStart, $goBackStart))
}

# 2. Fix the printed letter date: "11 mai 2012" -> "14 mai 2012"
$d.Content.Find.Execute("11 mai 2012", $true, $false, $false, $false, $false,
                         $true, 1, $false, "14 mai 2012", 2)
